$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Slot 5" (E2) changes from "Ginger Ale" to "Zuckersirup"
$ws.Range("E2").Value = "Zuckersirup"

# The column layout was re-fitted to the new content (widths no longer
# marked as auto "best fit" - explicit custom widths instead). Columns
# E (5) and J (10) are unaffected.
$ws.Columns.Item(1).ColumnWidth = 12.333333333333334
$ws.Columns.Item(2).ColumnWidth = 8
$ws.Columns.Item(3).ColumnWidth = 8.166666666666666
$ws.Columns.Item(4).ColumnWidth = 12
$ws.Columns.Item(6).ColumnWidth = 12.833333333333334
$ws.Columns.Item(7).ColumnWidth = 13
$ws.Columns.Item(8).ColumnWidth = 11
$ws.Columns.Item(9).ColumnWidth = 12.833333333333334

# Move/restore the active selection to E3
$ws.Range("E3").Select() | Out-Null
